# Case and Fatality Demographics Data Updated
# Monthly refresh of the three "Fatalities by ..." tabs (Age / Sex / RaceEth)
# with revised counts for recent months.
$wb = $excel.ActiveWorkbook

# --- Fatalities by Age ---
$ws = $wb.Worksheets.Item("Fatalities by Age")
$ws.Range("I9").Value2 = 1394
$ws.Range("L9").Value2 = 6674
$ws.Range("K14").Value2 = 2601
$ws.Range("L14").Value2 = 7036
$ws.Range("I15").Value2 = 6566
$ws.Range("K15").Value2 = 10777
$ws.Range("K18").Value2 = 1515
$ws.Range("L18").Value2 = 5242
$ws.Range("H20").Value2 = 229
$ws.Range("J20").Value2 = 267
$ws.Range("K20").Value2 = 257
$ws.Range("L20").Value2 = 1244
$ws.Range("J21").Value2 = 252
$ws.Range("L21").Value2 = 1010
$ws.Range("K22").Value2 = 110
$ws.Range("L22").Value2 = 620
$ws.Range("E24").Value2 = 129
$ws.Range("F24").Value2 = 404
$ws.Range("G24").Value2 = 849
$ws.Range("H24").Value2 = 1296
$ws.Range("I24").Value2 = 1490
$ws.Range("J24").Value2 = 1418
$ws.Range("K24").Value2 = 1270
$ws.Range("L24").Value2 = 6876
$ws.Range("E25").Value2 = 150
$ws.Range("F25").Value2 = 504
$ws.Range("G25").Value2 = 1039
$ws.Range("H25").Value2 = 1530
$ws.Range("I25").Value2 = 2007
$ws.Range("J25").Value2 = 1664
$ws.Range("L25").Value2 = 8288
$ws.Range("E26").Value2 = 82
$ws.Range("G26").Value2 = 430
$ws.Range("H26").Value2 = 750
$ws.Range("I26").Value2 = 1087
$ws.Range("J26").Value2 = 849
$ws.Range("K26").Value2 = 639
$ws.Range("L26").Value2 = 4084
$ws.Range("D27").Value2 = 2
$ws.Range("F27").Value2 = 72
$ws.Range("G27").Value2 = 161
$ws.Range("H27").Value2 = 309
$ws.Range("I27").Value2 = 398
$ws.Range("J27").Value2 = 354
$ws.Range("K27").Value2 = 320
$ws.Range("L27").Value2 = 1643
$ws.Range("D28").Value2 = 3
$ws.Range("E28").Value2 = 27
$ws.Range("F28").Value2 = 64
$ws.Range("G28").Value2 = 166
$ws.Range("H28").Value2 = 283
$ws.Range("I28").Value2 = 421
$ws.Range("J28").Value2 = 428
$ws.Range("K28").Value2 = 439
$ws.Range("L28").Value2 = 1834
$ws.Range("D29").Value2 = 60
$ws.Range("E29").Value2 = 549
$ws.Range("F29").Value2 = 1735
$ws.Range("G29").Value2 = 3858
$ws.Range("H29").Value2 = 6765
$ws.Range("I29").Value2 = 10305
$ws.Range("J29").Value2 = 10510
$ws.Range("K29").Value2 = 10482
$ws.Range("L29").Value2 = 44291
$ws.Range("D31").Value2 = 85
$ws.Range("E31").Value2 = 712
$ws.Range("F31").Value2 = 2247
$ws.Range("G31").Value2 = 5433
$ws.Range("H31").Value2 = 10281
$ws.Range("I31").Value2 = 16871
$ws.Range("J31").Value2 = 18676
$ws.Range("K31").Value2 = 21259
$ws.Range("L31").Value2 = 75603

# --- Fatalities by Sex ---
$ws = $wb.Worksheets.Item("Fatalities by Sex")
$ws.Range("C9").Value2 = 3943
$ws.Range("E9").Value2 = 6674
$ws.Range("C14").Value2 = 4025
$ws.Range("E14").Value2 = 7036
$ws.Range("B18").Value2 = 2085
$ws.Range("E18").Value2 = 5242
$ws.Range("C20").Value2 = 692
$ws.Range("E20").Value2 = 1244
$ws.Range("B21").Value2 = 427
$ws.Range("E21").Value2 = 1010
$ws.Range("C22").Value2 = 349
$ws.Range("E22").Value2 = 620
$ws.Range("B24").Value2 = 2803
$ws.Range("C24").Value2 = 4073
$ws.Range("E24").Value2 = 6876
$ws.Range("B25").Value2 = 3436
$ws.Range("C25").Value2 = 4852
$ws.Range("E25").Value2 = 8288
$ws.Range("B26").Value2 = 1786
$ws.Range("C26").Value2 = 2298
$ws.Range("E26").Value2 = 4084
$ws.Range("B27").Value2 = 685
$ws.Range("C27").Value2 = 958
$ws.Range("E27").Value2 = 1643
$ws.Range("B28").Value2 = 769
$ws.Range("C28").Value2 = 1065
$ws.Range("E28").Value2 = 1834
$ws.Range("B29").Value2 = 18482
$ws.Range("C29").Value2 = 25809
$ws.Range("E29").Value2 = 44291
$ws.Range("B31").Value2 = 31673
$ws.Range("C31").Value2 = 43929
$ws.Range("E31").Value2 = 75603

# --- Fatalities by RaceEth ---
$ws = $wb.Worksheets.Item("Fatalities by RaceEth")
$ws.Range("C9").Value2 = 660
$ws.Range("H9").Value2 = 6674
$ws.Range("F14").Value2 = 3590
$ws.Range("H14").Value2 = 7036
$ws.Range("C15").Value2 = 3177
$ws.Range("F15").Value2 = 11891
$ws.Range("D18").Value2 = 2380
$ws.Range("H18").Value2 = 5242
$ws.Range("C20").Value2 = 167
$ws.Range("F20").Value2 = 449
$ws.Range("H20").Value2 = 1244
$ws.Range("F21").Value2 = 385
$ws.Range("H21").Value2 = 1010
$ws.Range("F22").Value2 = 252
$ws.Range("H22").Value2 = 620
$ws.Range("C24").Value2 = 813
$ws.Range("D24").Value2 = 2539
$ws.Range("F24").Value2 = 3401
$ws.Range("G24").Value2 = 10
$ws.Range("H24").Value2 = 6876
$ws.Range("C25").Value2 = 1005
$ws.Range("D25").Value2 = 2976
$ws.Range("E25").Value2 = 60
$ws.Range("F25").Value2 = 4143
$ws.Range("H25").Value2 = 8288
$ws.Range("C26").Value2 = 477
$ws.Range("D26").Value2 = 1538
$ws.Range("F26").Value2 = 1973
$ws.Range("H26").Value2 = 4084
$ws.Range("B27").Value2 = 22
$ws.Range("C27").Value2 = 139
$ws.Range("D27").Value2 = 668
$ws.Range("F27").Value2 = 800
$ws.Range("H27").Value2 = 1643
$ws.Range("B28").Value2 = 27
$ws.Range("C28").Value2 = 141
$ws.Range("D28").Value2 = 705
$ws.Range("E28").Value2 = 21
$ws.Range("F28").Value2 = 939
$ws.Range("G28").Value2 = 1
$ws.Range("H28").Value2 = 1834
$ws.Range("B29").Value2 = 805
$ws.Range("C29").Value2 = 4840
$ws.Range("D29").Value2 = 17381
$ws.Range("E29").Value2 = 292
$ws.Range("F29").Value2 = 20944
$ws.Range("G29").Value2 = 29
$ws.Range("H29").Value2 = 44291
$ws.Range("B31").Value2 = 1382
$ws.Range("C31").Value2 = 8017
$ws.Range("D31").Value2 = 32859
$ws.Range("E31").Value2 = 463
$ws.Range("F31").Value2 = 32835
$ws.Range("G31").Value2 = 47
$ws.Range("H31").Value2 = 75603

# Leave "Fatalities by Age" as the active/selected sheet (first tab),
# clearing the stray B17 selection left over from the last edit session.
$wsActive = $wb.Worksheets.Item("Fatalities by Age")
$wsActive.Activate()
[void]$wsActive.Range("A1").Select()
